$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Price" (column D) values. Several of these are plain numeric-looking
# strings (e.g. "1.00", "57.95") that Excel would otherwise auto-convert to
# real numbers on assignment, so the cell is forced to Text format first and
# the format is put back to the default "Normal" style afterwards so no extra
# styling is left behind on the cell.
$priceUpdates = @{
    2 = "37.696.47"
    3 = "2.078.16"
    4 = "1.00"
    5 = "232.56"
    8 = "57.95"
    10 = "0.0780"
    12 = "14.88"
    13 = "2.386.65"
    14 = "21.21"
    15 = "0.765"
    17 = "2.077.49"
    18 = "37.669.80"
    19 = "6.15"
    20 = "70.25"
    21 = "0.0₃0830"
    22 = "227.74"
    26 = "9.91"
    27 = "169.58"
    28 = "0.131"
    29 = "19.35"
    35 = "2.52"
    36 = "1.82"
    37 = "3.32"
    39 = "5.34"
    41 = "98.49"
    42 = "0.0959"
    44 = "1.491.31"
    46 = "16.95"
    47 = "4.10"
    50 = "2.97"
    51 = "2.270.46"
}

foreach ($row in $priceUpdates.Keys) {
    $cell = $ws.Cells.Item($row, 4)
    $cell.NumberFormat = "@"
    $cell.Value = $priceUpdates[$row]
    $cell.Style = "Normal"
}

# New "Volume(1h)" (column E) values - plain text, never ambiguous as numbers.
$volumeUpdates = @{
    2 = "  -0.30%  "
    3 = "  -0.28%  "
    4 = "  -0.08%  "
    5 = "  -0.46%  "
    6 = "  -0.30%  "
    7 = "  +0.04%  "
    8 = "  -1.73%  "
    9 = "  -1.08%  "
    10 = "  -1.23%  "
    11 = "  +0.46%  "
    12 = "  +0.60%  "
    13 = "  -0.10%  "
    14 = "  -0.16%  "
    15 = "  -1.31%  "
    16 = "  +0.04%  "
    17 = "  +0.34%  "
    18 = "  -0.14%  "
    19 = "  -0.02%  "
    20 = "  -2.18%  "
    21 = "  -2.05%  "
    22 = "  -0.29%  "
    23 = "  +0.03%  "
    24 = "  +0.41%  "
    25 = "  -2.39%  "
    26 = "  +2.90%  "
    27 = "  -0.89%  "
    28 = "  -3.87%  "
    29 = "  -1.15%  "
    30 = "  -2.84%  "
    31 = "  -0.03%  "
    32 = "  -2.90%  "
    33 = "  -0.88%  "
    34 = "  -0.49%  "
    35 = "  +1.13%  "
    36 = "  -0.48%  "
    38 = "  +0.14%  "
    39 = "  -1.83%  "
    40 = "  +3.68%  "
    42 = "  -2.33%  "
    43 = "  +0.31%  "
    44 = "  +2.73%  "
    45 = "  +3.23%  "
    46 = "  -2.48%  "
    47 = "  -1.36%  "
    48 = "  -1.90%  "
    49 = "  -1.16%  "
    50 = "  -0.75%  "
    51 = "  -0.24%  "
}

foreach ($row in $volumeUpdates.Keys) {
    $ws.Cells.Item($row, 5).Value = $volumeUpdates[$row]
}

